$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns with latest scraped cryptos data.
# D-column values are forced to Text (quote-prefix cleared afterwards) so that
# numeric-looking strings like "1.000" / "215.02" are not auto-coerced to numbers
# by Excel, matching the original inlineStr/text cell representation.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '26.129.64'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +0.00%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.645.19'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -1.29%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.000'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.32%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '215.02'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +2.71%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.5227'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.46%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.000'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.32%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2609'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -0.24%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06327'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.18%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '20.86'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -1.07%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07694'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +2.22%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.630.71'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -2.38%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.425'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +0.03%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '1.861.62'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -1.77%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.5585'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +1.96%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.0₅8199'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +3.33%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '65.22'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -1.59%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '26.115.57'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -0.10%  '

$ws.Range("E19").Value = '  -0.12%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '4.750'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +0.62%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '189.19'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +1.60%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '10.23'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.37%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '6.191'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.33%  '

$ws.Range("E24").Value = '  -0.27%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '145.77'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -2.15%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '7.434'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.57%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.1208'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -2.98%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '15.89'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +0.42%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.396'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +3.30%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.05885'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -7.47%  '

$ws.Range("E31").Value = '  -1.16%  '

$ws.Range("E32").Value = '  -1.63%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.415'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +0.27%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.656'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +1.02%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.9872'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -1.52%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.763'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +0.52%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.391'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -0.67%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.5672'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -5.58%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01618'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +0.28%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.8575'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -1.02%  '

$ws.Range("E41").Value = '  -0.27%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '5.729'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -6.02%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.031.43'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -7.03%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '100.20'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +0.17%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.789.45'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -1.61%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0₈106'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -4.61%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '56.05'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +1.51%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.005'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +0.75%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '8.077'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +0.71%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.05154'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -1.46%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.4220'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -0.66%  '
